# Update countries & provincias Spain
# - Swap ranking of Singapur/Suiza (Singapur overtook Suiza in total cases)
# - Swap ranking of Ghana/Finlandia (Ghana overtook Finlandia in total cases)
# - Refresh Portugal's active/recovered numbers
# - Bump the "last updated" timestamp from 09:05 to 09:35

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / timestamp row
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 09:35"

# Row 28 now becomes Singapur (previously Suiza), row 29 now becomes Suiza
# (previously Singapur) - Singapur's case count passed Suiza's.
$ws.Range("A28:H28").Value = @(@("Singapur", 31068, 642, 12995, 18050, 0, 0, 23))
$ws.Range("A29:H29").Value = @(@("Suiza", 30707, 0, 27900, 904, 0, 0, 1903))

# Portugal (row 35) keeps its rank, only active/recovered figures refresh.
$ws.Range("D35").Value = 8977
$ws.Range("E35").Value = 10660

# Row 64 now becomes Ghana (previously Finlandia), row 65 now becomes
# Finlandia (previously Ghana) - Ghana's case count passed Finlandia's.
$ws.Range("A64:H64").Value = @(@("Ghana", 6617, 131, 1976, 4610, 0, 0, 31))
$ws.Range("A65:H65").Value = @(@("Finlandia", 6537, 0, 4800, 1431, 0, 0, 306))
